# Ocean freight price list update.
# The workbook holds a POL/POD/vehicle-type price table on Sheet1 (A:D,
# header in row 1). This script bumps the per-lane PRICE (column D) for a
# handful of rows, matching the data refresh captured in the source diff.
# (Window/theme/rupBuild/column-bestfit churn in the diff is just Excel
# re-saving the same file on a different machine and isn't a deliberate
# content edit, so it's intentionally left alone here.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value  = 1200   # New York  -> Rotterdam, SUV       1000 -> 1200
$ws.Range("D4").Value  = 1300   # New York  -> Rotterdam, LARGE SUV 1200 -> 1300
$ws.Range("D8").Value  = 1200   # Savannah  -> Rotterdam, SUV       1100 -> 1200
$ws.Range("D9").Value  = 1300   # Savannah  -> Rotterdam, LARGE SUV 1200 -> 1300
$ws.Range("D18").Value = 1500   # Houston   -> Rotterdam, SUV       1400 -> 1500
$ws.Range("D19").Value = 1600   # Houston   -> Rotterdam, LARGE SUV 1500 -> 1600
$ws.Range("D20").Value = 1700   # Houston   -> Rotterdam, PICKUP    1600 -> 1700
$ws.Range("D23").Value = 1500   # Indianapolis -> Rotterdam, SUV       1400 -> 1500
$ws.Range("D24").Value = 1600   # Indianapolis -> Rotterdam, LARGE SUV 1500 -> 1600
$ws.Range("D25").Value = 1700   # Indianapolis -> Rotterdam, PICKUP    1600 -> 1700
$ws.Range("D42").Value = 1700   # Savannah  -> Varna, CAR           1600 -> 1700
$ws.Range("D47").Value = 1700   # Miami     -> Varna, CAR           1600 -> 1700
$ws.Range("D52").Value = 2000   # Houston   -> Varna, CAR           1900 -> 2000
$ws.Range("D53").Value = 2100   # Houston   -> Varna, SUV           2000 -> 2100
$ws.Range("D54").Value = 2200   # Houston   -> Varna, LARGE SUV     2100 -> 2200
$ws.Range("D56").Value = 1000   # Houston   -> Varna, MOTORCYCLE     900 -> 1000

# Mirror the cursor position captured in the saved file.
$ws.Range("D62").Select()
